# Update the "Latest HO Xliff Generate Date" for the af3f0b90 file on the
# Overview sheet, and the corresponding handoff/handback timestamps for
# that same file on the zh-cn and de-de report sheets. This reflects a
# newly generated handback report for af3f0b90-2ce4-4fe3-aa15-a8a3920e9e6d.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-30 19:00:58"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-30 19:00:54"
$zhcn.Range("K3").Value = "2016-08-30 19:01:23"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-30 19:00:58"
$dede.Range("K3").Value = "2016-08-30 19:01:31"
